# Auto-generated edit script applying the diff changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('展览')
$ws.Range("F4").Value = 1658
$ws.Range("F5").Value = 557
$ws.Range("F6").Value = 110
$ws.Range("F7").Value = 341
$ws.Range("F8").Value = 506
$ws.Range("F9").Value = 240
$ws.Range("F10").Value = 1100
$ws.Range("F11").Value = 823
$ws.Range("F12").Value = 400
$ws.Range("F14").Value = 559
$ws.Range("F17").Value = 531
$ws.Range("F18").Value = 13
$ws.Range("F19").Value = 1259
$ws.Range("F21").Value = 1179
$ws.Range("F22").Value = 135
$ws.Range("F23").Value = 105
$ws.Range("F24").Value = 2446
$ws.Range("F25").Value = 235
$ws.Range("F26").Value = 113
$ws.Range("F32").Value = 175

$ws = $wb.Worksheets.Item('演出')
$ws.Range("F8").Value = 53
$ws.Range("F11").Value = 608
$ws.Range("F12").Value = 301
$ws.Range("F16").Value = 4
$ws.Range("F17").Value = 210
$ws.Range("F20").Value = 204
$ws.Range("F22").Value = 117
$ws.Range("F24").Value = 168
$ws.Range("F26").Value = 7

$ws = $wb.Worksheets.Item('本地生活')
$ws.Range("F2").Value = 1717
$ws.Range("F4").Value = 36
$ws.Range("F5").Value = 1959
$ws.Range("F6").Value = 2045
$ws.Range("F8").Value = 796

$ws = $wb.Worksheets.Item('全部类型')
$ws.Range("F2").Value = 1717
$ws.Range("F4").Value = 36
$ws.Range("F6").Value = 1959
$ws.Range("F7").Value = 2045
$ws.Range("F10").Value = 1658
$ws.Range("F13").Value = 796
$ws.Range("F15").Value = 557
$ws.Range("F18").Value = 110
$ws.Range("F19").Value = 53
$ws.Range("F20").Value = 341
$ws.Range("F21").Value = 506
$ws.Range("F22").Value = 240
$ws.Range("F24").Value = 1100
$ws.Range("F25").Value = 823
$ws.Range("F27").Value = 400
$ws.Range("F29").Value = 608
$ws.Range("F30").Value = 531
$ws.Range("F31").Value = 1259
$ws.Range("F32").Value = 301
$ws.Range("F34").Value = 1179
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = '2024-02-17'
$ws.Range("B35").NumberFormat = "General"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = '上海·少女番only2.0'
$ws.Range("C35").NumberFormat = "General"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '营口路699号(黄兴公园地铁站2号口旁) 花嫁丽舍'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '2024.02.17 10:00-02.17 17:00'
$ws.Range("E35").NumberFormat = "General"
$ws.Range("F35").Value = 135
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '80'
$ws.Range("G35").NumberFormat = "General"
$ws.Range("I35").NumberFormat = "@"
$ws.Range("I35").Value = 'https://show.bilibili.com/platform/detail.html?id=81148'
$ws.Range("I35").NumberFormat = "General"
$ws.Range("J35").NumberFormat = "@"
$ws.Range("J35").Value = '//i1.hdslb.com/bfs/openplatform/202401/j6eEZ18S1705657346664.jpeg'
$ws.Range("J35").NumberFormat = "General"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = '上海·SISP动漫游戏嘉年华'
$ws.Range("C36").NumberFormat = "General"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '年家浜路518号 周浦万达广场'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '2024.02.24 13:00-02.25 19:00'
$ws.Range("E36").NumberFormat = "General"
$ws.Range("F36").Value = 105
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '48'
$ws.Range("G36").NumberFormat = "General"
$ws.Range("H36").Value = $false
$ws.Range("I36").NumberFormat = "@"
$ws.Range("I36").Value = 'https://show.bilibili.com/platform/detail.html?id=80339'
$ws.Range("I36").NumberFormat = "General"
$ws.Range("J36").NumberFormat = "@"
$ws.Range("J36").Value = '//i0.hdslb.com/bfs/openplatform/202312/a8iuOufB1703832570508.jpeg'
$ws.Range("J36").NumberFormat = "General"
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = '2024-02-24'
$ws.Range("B37").NumberFormat = "General"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = '上海·第三届燃梦BACG PRO游戏动漫展-原X铁X崩同好交流'
$ws.Range("C37").NumberFormat = "General"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '盈浦街道淀山浦社区淀山湖大道851号青浦万达茂F3 万达汽车乐园(青浦万达茂店)'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '2024.02.24 11:00-02.25 16:30'
$ws.Range("E37").NumberFormat = "General"
$ws.Range("F37").Value = 2446
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '65'
$ws.Range("G37").NumberFormat = "General"
$ws.Range("H37").Value = $true
$ws.Range("I37").NumberFormat = "@"
$ws.Range("I37").Value = 'https://show.bilibili.com/platform/detail.html?id=77754'
$ws.Range("I37").NumberFormat = "General"
$ws.Range("J37").NumberFormat = "@"
$ws.Range("J37").Value = '//i1.hdslb.com/bfs/openplatform/202312/7eGZETK91701943985671.jpeg'
$ws.Range("J37").NumberFormat = "General"
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = '2024-02-25'
$ws.Range("B38").NumberFormat = "General"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = '上海·青山吉能见面会'
$ws.Range("C38").NumberFormat = "General"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '虹许路731号4号楼 THE BOXX•城市乐园'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '2024.02.25 14:30-02.25 19:30'
$ws.Range("E38").NumberFormat = "General"
$ws.Range("F38").Value = 194
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '380'
$ws.Range("G38").NumberFormat = "General"
$ws.Range("I38").NumberFormat = "@"
$ws.Range("I38").Value = 'https://show.bilibili.com/platform/detail.html?id=80142'
$ws.Range("I38").NumberFormat = "General"
$ws.Range("J38").NumberFormat = "@"
$ws.Range("J38").Value = '//i0.hdslb.com/bfs/openplatform/202312/1npuHFBM1703231674558.jpeg'
$ws.Range("J38").NumberFormat = "General"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = '上海·原神X星穹铁道ONLY'
$ws.Range("C39").NumberFormat = "General"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '逸仙路301号靠纪念路路口 上海宝丰联大酒店'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '2024.03.02 10:00-03.02 17:00'
$ws.Range("E39").NumberFormat = "General"
$ws.Range("F39").Value = 235
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '60'
$ws.Range("G39").NumberFormat = "General"
$ws.Range("H39").Value = $false
$ws.Range("I39").NumberFormat = "@"
$ws.Range("I39").Value = 'https://show.bilibili.com/platform/detail.html?id=80299'
$ws.Range("I39").NumberFormat = "General"
$ws.Range("J39").NumberFormat = "@"
$ws.Range("J39").Value = '//i2.hdslb.com/bfs/openplatform/202312/V0xu26Cl1703753850690.jpeg'
$ws.Range("J39").NumberFormat = "General"
$ws.Range("F40").Value = 210
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = '2024-03-02'
$ws.Range("B41").NumberFormat = "General"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = '上海·小山百代2024上海粉丝见面会'
$ws.Range("C41").NumberFormat = "General"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '宜昌路179号 万代南梦宫上海文化中心'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2024.03.02 13:00-03.02 20:00'
$ws.Range("E41").NumberFormat = "General"
$ws.Range("F41").Value = 210
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '380'
$ws.Range("G41").NumberFormat = "General"
$ws.Range("H41").Value = $true
$ws.Range("I41").NumberFormat = "@"
$ws.Range("I41").Value = 'https://show.bilibili.com/platform/detail.html?id=80924'
$ws.Range("I41").NumberFormat = "General"
$ws.Range("J41").NumberFormat = "@"
$ws.Range("J41").Value = '//i1.hdslb.com/bfs/openplatform/202401/FpA9OkKy1705467080070.jpeg'
$ws.Range("J41").NumberFormat = "General"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = '2024-03-03'
$ws.Range("B42").NumberFormat = "General"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = '上海·怀旧番ONLY'
$ws.Range("C42").NumberFormat = "General"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '逸仙路270号  上海宝丰联大酒店'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '2024.03.03 10:00-03.03 17:00'
$ws.Range("E42").NumberFormat = "General"
$ws.Range("F42").Value = 113
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '60'
$ws.Range("G42").NumberFormat = "General"
$ws.Range("H42").Value = $false
$ws.Range("I42").NumberFormat = "@"
$ws.Range("I42").Value = 'https://show.bilibili.com/platform/detail.html?id=80575'
$ws.Range("I42").NumberFormat = "General"
$ws.Range("J42").NumberFormat = "@"
$ws.Range("J42").Value = '//i1.hdslb.com/bfs/openplatform/202401/y4uWdyPT1704700763902.jpeg'
$ws.Range("J42").NumberFormat = "General"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = '上海·第五十三届燃梦星辰动漫嘉年华-随机宅舞'
$ws.Range("C43").NumberFormat = "General"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '周家嘴路3608号 宝龙旭辉广场'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '2024.03.09 10:20-03.10 16:30'
$ws.Range("E43").NumberFormat = "General"
$ws.Range("F43").Value = 37
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '58'
$ws.Range("G43").NumberFormat = "General"
$ws.Range("H43").Value = $true
$ws.Range("I43").NumberFormat = "@"
$ws.Range("I43").Value = 'https://show.bilibili.com/platform/detail.html?id=80571'
$ws.Range("I43").NumberFormat = "General"
$ws.Range("J43").NumberFormat = "@"
$ws.Range("J43").Value = '//i0.hdslb.com/bfs/openplatform/202401/SHH70VXN1704700240858.jpeg'
$ws.Range("J43").NumberFormat = "General"
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = '2024-03-09'
$ws.Range("B44").NumberFormat = "General"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = '上海·青山刚昌ONLY【名侦探柯南&魔术快斗】'
$ws.Range("C44").NumberFormat = "General"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '漕宝路1688号 诺宝中心酒店'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '2024.03.09 10:00-03.09 16:30'
$ws.Range("E44").NumberFormat = "General"
$ws.Range("F44").Value = 853
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '73'
$ws.Range("G44").NumberFormat = "General"
$ws.Range("I44").NumberFormat = "@"
$ws.Range("I44").Value = 'https://show.bilibili.com/platform/detail.html?id=76410'
$ws.Range("I44").NumberFormat = "General"
$ws.Range("J44").NumberFormat = "@"
$ws.Range("J44").Value = '//i2.hdslb.com/bfs/openplatform/202309/fVXMrcHy1693971682397.jpeg'
$ws.Range("J44").NumberFormat = "General"
$ws.Range("F45").Value = 117
$ws.Range("F46").Value = 117
$ws.Range("F49").Value = 168
$ws.Range("F51").Value = 175
